# Update DateBase/orders/Fresh bloom Flowers_2025-11-20.xlsx
# - Orders sheet: F41 number changes 1 -> 15, and 20 new order rows (42-61)
#   are appended (PackageID in col A, FlowerName in col C, Number in col F).
# - Summary sheet: G2 (TotalNumber) string gets the new numbers appended.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Orders")

# Helper: write a value as genuine TEXT (matches the workbook's existing
# "numberStoredAsText" convention for every cell in this sheet) instead of
# letting Excel auto-detect numeric-looking strings as numbers.
function Set-TextValue($sheet, [int]$row, [int]$col, [string]$val) {
    $cell = $sheet.Cells.Item($row, $col)
    if ($val -match '^-?[0-9]+(\.[0-9]+)?$') {
        $cell.Value = "'" + $val
    } else {
        $cell.Value = $val
    }
}

# Existing row 41: Number 1 -> 15
Set-TextValue $ws 41 6 "15"

# New rows 42-61
$rows = @(
    @{ R = 42; A = "8";  C = "586_洋牡丹白_undefined_undefined_1bunch";            F = "15" },
    @{ R = 43; A = "";   C = "590_洋牡丹粉_undefined_undefined_1bunch";            F = "10" },
    @{ R = 44; A = "";   C = "585_洋牡丹红_undefined_undefined_1bunch";            F = "5"  },
    @{ R = 45; A = "";   C = "649_洋牡丹樱花粉_undefined_undefined_1bunch";         F = "10" },
    @{ R = 46; A = "";   C = "648_洋牡丹河内_undefined_undefined_1bunch";          F = "10" },
    @{ R = 47; A = "";   C = "480_蝴蝶洋牡丹红_butterfly  Ranunculus_undefined_1bunch"; F = "10" },
    @{ R = 48; A = "";   C = "419_松虫草红_scabiosa watermelon_undefined_1bunch";  F = "15" },
    @{ R = 49; A = "9";  C = "721_银扇干花_undefined_undefined_1bunch";            F = "40" },
    @{ R = 50; A = "10"; C = "512_松虫草粉_scabiosa pink_undefined_1bunch";        F = "9"  },
    @{ R = 51; A = "";   C = "418_松虫草白_scabiosa white_undefined_1bunch";       F = "10" },
    @{ R = 52; A = "";   C = "514_松虫草紫_scabiosa purple_undefined_1bunch";      F = "10" },
    @{ R = 53; A = "";   C = "395_豌豆花 混色_sweetpea mix colors_undefined_1bunch"; F = "26" },
    @{ R = 54; A = "11"; C = "373_龙柳_Salix`n_undefined_1bunch";                  F = "10" },
    @{ R = 55; A = "12"; C = "373_龙柳_Salix`n_undefined_1bunch";                  F = "10" },
    @{ R = 56; A = "1";  C = "137_凯瑟琳_Catherine_Rosa rugosa Thunb._20stems";    F = "13" },
    @{ R = 57; A = "";   C = "135_甜蜜曼塔_sweet menta_Rosa rugosa Thunb._20stems"; F = "17" },
    @{ R = 58; A = "";   C = "152_白荔枝_White Ohara_Rosa rugosa Thunb._20stems";  F = "8"  },
    @{ R = 59; A = "";   C = "412_紫罗兰粉_violet pink_undefined_1bunch";          F = "10" },
    @{ R = 60; A = "";   C = "600_康乃馨复古红_vintage red_undefined_20stems";     F = "5"  },
    @{ R = 61; A = "2";  C = "";                                                   F = ""   }
)

foreach ($row in $rows) {
    if ($row.A -ne "") { Set-TextValue $ws $row.R 1 $row.A }
    if ($row.C -ne "") { Set-TextValue $ws $row.R 3 $row.C }
    if ($row.F -ne "") { Set-TextValue $ws $row.R 6 $row.F }
}

# Summary sheet: append the new Number values onto the TotalNumber string
$ws2 = $wb.Worksheets.Item("Summary")
$g2 = $ws2.Range("G2").Value2
$suffix = "515105101010154091010261010131781050"
$newG2 = [string]::Concat($g2, $suffix)
Set-TextValue $ws2 2 7 $newG2
